# "Generate Report for Archive"
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview summary columns for
#   zh-cn/de-de, and the per-language "Status" column on each language
#   sheet).
# - Shrink the "Status"-related columns that held the old, longer text so
#   they fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: "Status" column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: "Status" column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5

Write-Output "Updated status text and column widths."
